# Update res_line/pl_mw values for the "case with 380 kV done" run.
# Columns B, D, E, F, G, J, L, M, O over rows 2-25 (A is the index column,
# C/H/I/K/N remain 0 and are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (index 0)
$ws.Cells.Item(2, 2).Value = 1.532262534272604
$ws.Cells.Item(2, 4).Value = 0.2504141269181517
$ws.Cells.Item(2, 5).Value = 0.3166662929811181
$ws.Cells.Item(2, 6).Value = 1.396569297610434
$ws.Cells.Item(2, 7).Value = 0.002455192029909734
$ws.Cells.Item(2, 10).Value = 0.4623035817220682
$ws.Cells.Item(2, 12).Value = 0.5388787369807062
$ws.Cells.Item(2, 13).Value = 0.4418470755331612
$ws.Cells.Item(2, 15).Value = 3.367291297544455

# Row 3 (index 1)
$ws.Cells.Item(3, 2).Value = 1.446322473854764
$ws.Cells.Item(3, 4).Value = 0.2497068241671556
$ws.Cells.Item(3, 5).Value = 0.3097559937271157
$ws.Cells.Item(3, 6).Value = 1.415502031845392
$ws.Cells.Item(3, 7).Value = 0.002458610405048538
$ws.Cells.Item(3, 10).Value = 0.4461774775234346
$ws.Cells.Item(3, 12).Value = 0.4800428348249568
$ws.Cells.Item(3, 13).Value = 0.4063941292999118
$ws.Cells.Item(3, 15).Value = 3.393874959413012

# Row 4 (index 2)
$ws.Cells.Item(4, 2).Value = 1.393697605099419
$ws.Cells.Item(4, 4).Value = 0.2493303591073328
$ws.Cells.Item(4, 5).Value = 0.305585813266795
$ws.Cells.Item(4, 6).Value = 1.428229946002254
$ws.Cells.Item(4, 7).Value = 0.002460821781583182
$ws.Cells.Item(4, 10).Value = 0.4363728616422691
$ws.Cells.Item(4, 12).Value = 0.4437922635034113
$ws.Cells.Item(4, 13).Value = 0.3845983191157671
$ws.Cells.Item(4, 15).Value = 3.413019331086389

# Row 5 (index 3)
$ws.Cells.Item(5, 2).Value = 1.372289654584421
$ws.Cells.Item(5, 4).Value = 0.2491915867993697
$ws.Cells.Item(5, 5).Value = 0.3039049796312696
$ws.Cells.Item(5, 6).Value = 1.433693633013419
$ws.Cells.Item(5, 7).Value = 0.002461751306563035
$ws.Cells.Item(5, 10).Value = 0.4324022274982724
$ws.Cells.Item(5, 12).Value = 0.4289893019062561
$ws.Cells.Item(5, 13).Value = 0.3757099245626421
$ws.Cells.Item(5, 15).Value = 3.421528929062077

# Row 6 (index 4)
$ws.Cells.Item(6, 2).Value = 1.368737158894618
$ws.Cells.Item(6, 4).Value = 0.2491694309709658
$ws.Cells.Item(6, 5).Value = 0.3036270062030297
$ws.Cells.Item(6, 6).Value = 1.434617588935566
$ws.Cells.Item(6, 7).Value = 0.002461907369717932
$ws.Cells.Item(6, 10).Value = 0.4317444200309666
$ws.Cells.Item(6, 12).Value = 0.4265294623912723
$ws.Cells.Item(6, 13).Value = 0.3742336404217426
$ws.Cells.Item(6, 15).Value = 3.422984666136671

# Row 7 (index 5)
$ws.Cells.Item(7, 2).Value = 1.39340873808564
$ws.Cells.Item(7, 4).Value = 0.2493284281629755
$ws.Cells.Item(7, 5).Value = 0.3055630695193372
$ws.Cells.Item(7, 6).Value = 1.428302510383581
$ws.Cells.Item(7, 7).Value = 0.002460834202444519
$ws.Cells.Item(7, 10).Value = 0.4363192111370893
$ws.Cells.Item(7, 12).Value = 0.4435927481184763
$ws.Cells.Item(7, 13).Value = 0.3844784724566708
$ws.Cells.Item(7, 15).Value = 3.413131229361653

# Row 8 (index 6)
$ws.Cells.Item(8, 2).Value = 1.502601494906457
$ws.Cells.Item(8, 4).Value = 0.2501583074824367
$ws.Cells.Item(8, 5).Value = 0.3142686642769803
$ws.Cells.Item(8, 6).Value = 1.402868087715156
$ws.Cells.Item(8, 7).Value = 0.002456347392940507
$ws.Cells.Item(8, 10).Value = 0.4567234655260251
$ws.Cells.Item(8, 12).Value = 0.5186185944734518
$ws.Cells.Item(8, 13).Value = 0.4296289082746014
$ws.Cells.Item(8, 15).Value = 3.375870768113742

# Row 9 (index 7)
$ws.Cells.Item(9, 2).Value = 1.717818006750178
$ws.Cells.Item(9, 4).Value = 0.252240326692494
$ws.Cells.Item(9, 5).Value = 0.3319079147637751
$ws.Cells.Item(9, 6).Value = 1.361763694608598
$ws.Cells.Item(9, 7).Value = 0.002448437265073018
$ws.Cells.Item(9, 10).Value = 0.4974869114472398
$ws.Cells.Item(9, 12).Value = 0.6647199163314497
$ws.Cells.Item(9, 13).Value = 0.5179326913216045
$ws.Cells.Item(9, 15).Value = 3.325262445694193

# Row 10 (index 8)
$ws.Cells.Item(10, 2).Value = 1.876562890643584
$ws.Cells.Item(10, 4).Value = 0.2540419213365368
$ws.Cells.Item(10, 5).Value = 0.3452020509945726
$ws.Cells.Item(10, 6).Value = 1.336939125632384
$ws.Cells.Item(10, 7).Value = 0.002443161748828554
$ws.Cells.Item(10, 10).Value = 0.5278731767575948
$ws.Cells.Item(10, 12).Value = 0.771405117264095
$ws.Cells.Item(10, 13).Value = 0.5826491320586058
$ws.Cells.Item(10, 15).Value = 3.301870832204855

# Row 11 (index 9)
$ws.Cells.Item(11, 2).Value = 1.9489089988827
$ws.Cells.Item(11, 4).Value = 0.2549194947447404
$ws.Cells.Item(11, 5).Value = 0.3513202883429045
$ws.Cells.Item(11, 6).Value = 1.326818668159667
$ws.Cells.Item(11, 7).Value = 0.00244087698737773
$ws.Cells.Item(11, 10).Value = 0.5417876955578436
$ws.Cells.Item(11, 12).Value = 0.8197910338838312
$ws.Cells.Item(11, 13).Value = 0.6120526015048853
$ws.Cells.Item(11, 15).Value = 3.294245188324112

# Row 12 (index 10)
$ws.Cells.Item(12, 2).Value = 1.976322699006687
$ws.Cells.Item(12, 4).Value = 0.2552600587349474
$ws.Cells.Item(12, 5).Value = 0.3536470524305813
$ws.Cells.Item(12, 6).Value = 1.323155381142712
$ws.Cells.Item(12, 7).Value = 0.00244002826826281
$ws.Cells.Item(12, 10).Value = 0.5470695522045901
$ws.Cells.Item(12, 12).Value = 0.838091885600079
$ws.Cells.Item(12, 13).Value = 0.6231813263430013
$ws.Cells.Item(12, 15).Value = 3.29179283520071

# Row 13 (index 11)
$ws.Cells.Item(13, 2).Value = 1.970417895751837
$ws.Cells.Item(13, 4).Value = 0.2551863469771547
$ws.Cells.Item(13, 5).Value = 0.3531455047451502
$ws.Cells.Item(13, 6).Value = 1.323936805943589
$ws.Cells.Item(13, 7).Value = 0.002440210323953682
$ws.Cells.Item(13, 10).Value = 0.5459314507503166
$ws.Cells.Item(13, 12).Value = 0.8341514552863885
$ws.Cells.Item(13, 13).Value = 0.6207848190961158
$ws.Cells.Item(13, 15).Value = 3.292301605428207

# Row 14 (index 12)
$ws.Cells.Item(14, 2).Value = 1.951163990678083
$ws.Cells.Item(14, 4).Value = 0.254947348439913
$ws.Cells.Item(14, 5).Value = 0.3515115151748702
$ws.Cells.Item(14, 6).Value = 1.32651389509563
$ws.Cells.Item(14, 7).Value = 0.002440806833055253
$ws.Cells.Item(14, 10).Value = 0.5422219844917038
$ws.Cells.Item(14, 12).Value = 0.8212970988096231
$ws.Cells.Item(14, 13).Value = 0.6129682860891279
$ws.Cells.Item(14, 15).Value = 3.294034697641365

# Row 15 (index 13)
$ws.Cells.Item(15, 2).Value = 1.939372700019419
$ws.Cells.Item(15, 4).Value = 0.2548020259436328
$ws.Cells.Item(15, 5).Value = 0.3505119342086402
$ws.Cells.Item(15, 6).Value = 1.328114475393811
$ws.Cells.Item(15, 7).Value = 0.002441174354744979
$ws.Cells.Item(15, 10).Value = 0.5399514739754352
$ws.Cells.Item(15, 12).Value = 0.8134205633357681
$ws.Cells.Item(15, 13).Value = 0.6081796728720974
$ws.Cells.Item(15, 15).Value = 3.295153008014921

# Row 16 (index 14)
$ws.Cells.Item(16, 2).Value = 1.871837458959249
$ws.Cells.Item(16, 4).Value = 0.253985729016577
$ws.Cells.Item(16, 5).Value = 0.3448036110725852
$ws.Cells.Item(16, 6).Value = 1.337624157288488
$ws.Cells.Item(16, 7).Value = 0.002443313372602025
$ws.Cells.Item(16, 10).Value = 0.5269656359369037
$ws.Cells.Item(16, 12).Value = 0.7682399750583215
$ws.Cells.Item(16, 13).Value = 0.5807267713153692
$ws.Cells.Item(16, 15).Value = 3.302430008460249

# Row 17 (index 15)
$ws.Cells.Item(17, 2).Value = 1.830439776238393
$ws.Cells.Item(17, 4).Value = 0.2534997479122438
$ws.Cells.Item(17, 5).Value = 0.3413196663223772
$ws.Cells.Item(17, 6).Value = 1.343758684353048
$ws.Cells.Item(17, 7).Value = 0.002444655012205561
$ws.Cells.Item(17, 10).Value = 0.5190223972980732
$ws.Cells.Item(17, 12).Value = 0.740485186389094
$ws.Cells.Item(17, 13).Value = 0.563875631444219
$ws.Cells.Item(17, 15).Value = 3.307667728543748

# Row 18 (index 16)
$ws.Cells.Item(18, 2).Value = 1.806641444833417
$ws.Cells.Item(18, 4).Value = 0.2532256911862873
$ws.Cells.Item(18, 5).Value = 0.3393224592485211
$ws.Cells.Item(18, 6).Value = 1.347397419091088
$ws.Cells.Item(18, 7).Value = 0.002445437526232537
$ws.Cells.Item(18, 10).Value = 0.5144623097827008
$ws.Cells.Item(18, 12).Value = 0.7245077026308877
$ws.Cells.Item(18, 13).Value = 0.5541799080796466
$ws.Cells.Item(18, 15).Value = 3.310964040439103

# Row 19 (index 17)
$ws.Cells.Item(19, 2).Value = 1.798585926764588
$ws.Cells.Item(19, 4).Value = 0.2531338421838711
$ws.Cells.Item(19, 5).Value = 0.3386473915356731
$ws.Cells.Item(19, 6).Value = 1.348648363158681
$ws.Cells.Item(19, 7).Value = 0.002445704335980377
$ws.Cells.Item(19, 10).Value = 0.5129198425419759
$ws.Cells.Item(19, 12).Value = 0.7190956819001428
$ws.Cells.Item(19, 13).Value = 0.5508965353790245
$ws.Cells.Item(19, 15).Value = 3.312128790870275

# Row 20 (index 18)
$ws.Cells.Item(20, 2).Value = 1.834845343295399
$ws.Cells.Item(20, 4).Value = 0.2535509163491412
$ws.Cells.Item(20, 5).Value = 0.3416898500664303
$ws.Cells.Item(20, 6).Value = 1.343094232255304
$ws.Cells.Item(20, 7).Value = 0.00244451107118287
$ws.Cells.Item(20, 10).Value = 0.5198670763146822
$ws.Cells.Item(20, 12).Value = 0.743441151650245
$ws.Cells.Item(20, 13).Value = 0.5656698189231264
$ws.Cells.Item(20, 15).Value = 3.307080788541811

# Row 21 (index 19)
$ws.Cells.Item(21, 2).Value = 1.956818855867311
$ws.Cells.Item(21, 4).Value = 0.2550173250773042
$ws.Cells.Item(21, 5).Value = 0.3519911901890751
$ws.Cells.Item(21, 6).Value = 1.325752347742963
$ws.Cells.Item(21, 7).Value = 0.002440631177431879
$ws.Cells.Item(21, 10).Value = 0.5433112026915694
$ws.Cells.Item(21, 12).Value = 0.8250733347733217
$ws.Cells.Item(21, 13).Value = 0.6152643507927706
$ws.Cells.Item(21, 15).Value = 3.293513818548405

# Row 22 (index 20)
$ws.Cells.Item(22, 2).Value = 2.036638794213388
$ws.Cells.Item(22, 4).Value = 0.2560237276156698
$ws.Cells.Item(22, 5).Value = 0.3587814258094824
$ws.Cells.Item(22, 6).Value = 1.315404372835921
$ws.Cells.Item(22, 7).Value = 0.002438191407600666
$ws.Cells.Item(22, 10).Value = 0.5587072942765019
$ws.Cells.Item(22, 12).Value = 0.8782969565544363
$ws.Cells.Item(22, 13).Value = 0.6476436145211721
$ws.Cells.Item(22, 15).Value = 3.287185147475611

# Row 23 (index 21)
$ws.Cells.Item(23, 2).Value = 1.994028345513129
$ws.Cells.Item(23, 4).Value = 0.2554822296181527
$ws.Cells.Item(23, 5).Value = 0.3551521471782735
$ws.Cells.Item(23, 6).Value = 1.320836893901607
$ws.Cells.Item(23, 7).Value = 0.002439484804555363
$ws.Cells.Item(23, 10).Value = 0.5504834857521814
$ws.Cells.Item(23, 12).Value = 0.849902484842147
$ws.Cells.Item(23, 13).Value = 0.6303654324034653
$ws.Cells.Item(23, 15).Value = 3.290330074307434

# Row 24 (index 22)
$ws.Cells.Item(24, 2).Value = 1.832853580884603
$ws.Cells.Item(24, 4).Value = 0.2535277664572675
$ws.Cells.Item(24, 5).Value = 0.3415224720385623
$ws.Cells.Item(24, 6).Value = 1.34339428247722
$ws.Cells.Item(24, 7).Value = 0.002444576112199059
$ws.Cells.Item(24, 10).Value = 0.519485176382517
$ws.Cells.Item(24, 12).Value = 0.7421048247539375
$ws.Cells.Item(24, 13).Value = 0.5648586909565836
$ws.Cells.Item(24, 15).Value = 3.307345256229752

# Row 25 (index 23)
$ws.Cells.Item(25, 2).Value = 1.659483696224186
$ws.Cells.Item(25, 4).Value = 0.2516289580379976
$ws.Cells.Item(25, 5).Value = 0.3270764633293766
$ws.Cells.Item(25, 6).Value = 1.371941584648098
$ws.Cells.Item(25, 7).Value = 0.002450482625279462
$ws.Cells.Item(25, 10).Value = 0.4863810789566259
$ws.Cells.Item(25, 12).Value = 0.6253085642076996
$ws.Cells.Item(25, 13).Value = 0.4940711905268813
$ws.Cells.Item(25, 15).Value = 3.336538712667391
